# Generate Report for Handback
#
# The localization-status workbook tracks, per target language, the
# handoff/handback lifecycle of each source file. This run represents the
# files coming back from translation "in sync with en-US": the Status
# column flips from "Ready for handoff" to "Handed back: in sync with
# en-US", and each data row gains its "Latest Target File" / "Latest
# Handback File" / "Latest Handback DateTime" values.

$wb = $excel.ActiveWorkbook

$statusOld = "Ready for handoff"
$statusNew = "Handed back: in sync with en-US"

$mdFile  = "c0d4e954-9144-49ac-a031-20fb35067ae8.md"
$mdUrl   = "https://github.com/OpenLocalizationTest/oltest/blob/c88074a300af5880615ff9bbd1db62123137ca64/e2e/c0d4e954-9144-49ac-a031-20fb35067ae8.md"

$langs = @(
    @{
        Sheet = "zh-cn"
        XlfFile = "c0d4e954-9144-49ac-a031-20fb35067ae8.8028f1ebe34171eae2f288cc2c372639838fb653.zh-cn.xlf"
        XlfUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/849515db0e527731aeec36cc8f4b9d8bbd91ee5b/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/c0d4e954-9144-49ac-a031-20fb35067ae8.8028f1ebe34171eae2f288cc2c372639838fb653.zh-cn.xlf"
        HandbackDateTime = "2016-03-09 23:10:08"
    },
    @{
        Sheet = "de-de"
        XlfFile = "c0d4e954-9144-49ac-a031-20fb35067ae8.8028f1ebe34171eae2f288cc2c372639838fb653.de-de.xlf"
        XlfUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d229c610e58cbd0282aefcf8ef33f3d8f259200e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/c0d4e954-9144-49ac-a031-20fb35067ae8.8028f1ebe34171eae2f288cc2c372639838fb653.de-de.xlf"
        HandbackDateTime = "2016-03-09 23:10:20"
    }
)

foreach ($lang in $langs) {
    $ws = $wb.Worksheets.Item($lang.Sheet)

    # Status moves from "ready for handoff" to "handed back" for both data rows.
    $ws.Range("B2").Value = $statusNew
    $ws.Range("B3").Value = $statusNew

    # Row 2: Latest Target File / Latest Handback File / Latest Handback DateTime.
    $ws.Hyperlinks.Add($ws.Range("E2"), $mdUrl, "", "", $mdFile)
    $ws.Hyperlinks.Add($ws.Range("F2"), $lang.XlfUrl, "", "", $lang.XlfFile)
    $ws.Range("G2").Value = $lang.HandbackDateTime

    # Row 3 (the dependent file) reports back against the same handed-back
    # source/target pair as row 2.
    $ws.Hyperlinks.Add($ws.Range("E3"), $mdUrl, "", "", $mdFile)
    $ws.Hyperlinks.Add($ws.Range("F3"), $lang.XlfUrl, "", "", $lang.XlfFile)
    $ws.Range("G3").Value = $lang.HandbackDateTime
}

$wb.Save()
